$wb = $excel.ActiveWorkbook

# --- Update "Ready for handoff" -> "In Translation" everywhere it appears ---
# (Overview sheet: zh-cn/de-de status columns; per-language sheets: Status column)
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Narrow the status columns from ~17.22 to ~13.41 characters wide ---
$wsOverview.Range("E1:F1").ColumnWidth = 12.5
$wsZhCn.Range("C1").ColumnWidth = 12.5
$wsDeDe.Range("C1").ColumnWidth = 12.5
